# "made proper layer arch"
#
# Slide 24 ("The Course Project" - console output bullet):
#   Merge the "...to the " / "console" / "." runs so the sentence reads
#   as "...to the console" in one run, followed by a separate "." run.
#
# Slide 25 ("The Course Project" - next steps bullets):
#   Split "Make sure to make the program modular" into two runs, and
#   add a new bullet "We will add a database next session" right after
#   it (before "We will build an interface at the end of the course").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 24: "...to the console."
# ---------------------------------------------------------------------
$slide24 = $p.Slides.Item(24)
$shape24 = $slide24.Shapes.Item(2)
$tf24 = $shape24.TextFrame
$tr24 = $tf24.TextRange

$para = $tr24.Paragraphs(3, 1)
$sentence = "The program must then print the files in that folder to the console."

# Re-write the whole sentence (prefix + "console") as a single run while
# keeping the trailing period - this merges the old "to the " + "console"
# runs into one clean run.
$merged = $para.Replace($sentence, $sentence)

# Now split the trailing "." back out into its own run so it carries
# fresh (dirty="0" smtClean="0") run properties, matching the edit.
$tr24b = $tf24.TextRange
$para2 = $tr24b.Paragraphs(3, 1)
$lastChar = $para2.Characters($para2.Length - 1, 1)
$lastChar.Delete()

$tr24c = $tf24.TextRange
$para3 = $tr24c.Paragraphs(3, 1)
$null = $para3.InsertAfter(".")

# ---------------------------------------------------------------------
# Slide 25: modular bullet split + new "database" bullet
# ---------------------------------------------------------------------
$slide25 = $p.Slides.Item(25)
$shape25 = $slide25.Shapes.Item(2)
$tf25 = $shape25.TextFrame
$tr25 = $tf25.TextRange

$modularPara = $tr25.Paragraphs(2, 1)
$word = "modular"
$start = $modularPara.Length - $word.Length
$tail = $modularPara.Characters($start, $word.Length)
$tail.Delete()

$tr25b = $tf25.TextRange
$modularPara2 = $tr25b.Paragraphs(2, 1)
$null = $modularPara2.InsertAfter($word)

# Add the new bullet directly after the (now two-run) "modular" paragraph.
$tr25c = $tf25.TextRange
$modularPara3 = $tr25c.Paragraphs(2, 1)
$null = $modularPara3.InsertAfter("`rWe will add a database next session")
